$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O4:O8").ClearContents()
$ws.Range("O4:O8").Select()
Write-Host "Cleared O4:O8 and updated selection"
